$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty placeholder cells in row 2 (G2, H2, L2) so the cells
# are no longer present at all (matches diff removing those <c> nodes).
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("L2").ClearContents()

# Append a new data row (row 3) with the same shape as row 2.
$ws.Range("A3").Value = "Sumanth  Ratna"
$ws.Range("B3").Value = "21:30:59.464162"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "9999"
$ws.Range("D3").Value = "1A"
$ws.Range("E3").Value = "ku"
$ws.Range("F3").Value = "Augusta County"
$ws.Range("G3").Formula = '=""'
$ws.Range("H3").Formula = '=""'
$ws.Range("I3").Value = "7214 Bull Run Post Office Rd, Centreville, 20121"
$ws.Range("J3").Value = "127.0.0.1"
$ws.Range("K3").Value = "a9c7391c8a"
$ws.Range("L3").Formula = '=""'
